$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.507.71'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.35%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.922.75'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.08%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.40%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.79'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.12%  '

$ws.Range("E6").Value = '  +0.43%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4841'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4101'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08183'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.026'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.77'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +6.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.901.02'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.068'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.255'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.54%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.47'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.18%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06785'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.77%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.007'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.47%  '

$ws.Range("E18").Value = '  +1.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.84'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.74%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.007'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.548.38'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.639'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.77'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.44%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.186'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.19%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.178.55'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.726'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +10.86%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.13'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.40%  '

$ws.Range("E28").Value = '  +2.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.124'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.64'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.08%  '

$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09592'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.533'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.68%  '

$ws.Range("E34").Value = '  +0.75%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.394'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02290'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.92%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06156'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.59%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.182'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6000'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.56%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.062'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.73%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.83'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +8.35%  '

$ws.Range("E42").Value = '  +0.55%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1867'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.61%  '

$ws.Range("E44").Value = '  -1.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.279'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.67%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07608'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.56%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.46'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5607'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.76%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.965'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.66%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '117.34'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.54%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.443'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.12%  '
